# Stingmon database rows + view-state update
# (setting test scenario for fusion implementation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows for Stingmon skills -------------------------------------
# Row 64: ID 62 - NaturalSpectrum / PassiveSkill(Stingmon)
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = "NaturalSpectrum"
$ws.Range("C64").Value = "PassiveSkill(Stingmon)"
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0

# Row 65: ID 63 - SpikingStrike / DamageSkill(Stingmon)
$ws.Range("A65").Value = 63
$ws.Range("B65").Value = "SpikingStrike"
$ws.Range("C65").Value = "DamageSkill(Stingmon)"
$ws.Range("D65").Value = 45
$ws.Range("E65").Value = 3

# Match the "Bom" (Good) cell style used by the surrounding skill rows -
# copy the format from the row above (reuses the existing style, same as
# the rest of the table) instead of re-applying a named style (which would
# mint a duplicate style entry).
$ws.Range("A63:E63").Copy()
$ws.Range("A64:E65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Scroll / selection state as left by the author on save ----------------
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("C59").Select()
